# Generate Report for Handoff
#
# The two localized source files (681d4a9c-...md and 8681479b-...md) have
# finished handoff: status flips from "In Translation" to "Ready for
# handoff", new handoff timestamps/hashes are recorded, and the report rows
# are re-sorted so the 8681479b file now sorts first (row 2) and the
# 681d4a9c file sorts second (row 3) on every sheet.

$wb = $excel.ActiveWorkbook

# Helper: some literal values ("True"/"False") would otherwise be smart-typed
# as real booleans instead of the plain text the report stores them as. A
# leading apostrophe forces text entry (matching typing it into Excel by
# hand); re-applying the Normal style afterwards clears the resulting
# "quote prefix" cell flag so the cell format stays the same as every other
# plain-text cell.
function Set-TextValue($range, [string]$value) {
    $range.Value = "'" + $value
    $range.Style = "Normal"
}

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")

# Row 2 now describes the 8681479b file, row 3 the 681d4a9c file (swapped).
$ov.Range("A2").Value = "8681479b-ec15-4532-9f35-a51ee3b9248b.md"
$ov.Range("C2").Value = ".md"
$ov.Range("E2").Value = "Ready for handoff"
$ov.Range("F2").Value = "Ready for handoff"
$ov.Range("G2").Value = "2017-01-03 05:57:09"

$ov.Range("A3").Value = "681d4a9c-e7c4-4811-aa17-cc08b7db6c9f.md"
$ov.Range("C3").Value = ".md"
$ov.Range("E3").Value = "Ready for handoff"
$ov.Range("F3").Value = "Ready for handoff"
$ov.Range("G3").Value = "2017-01-03 05:57:09"

# Hyperlinks: the underlying relationship targets (rId2 -> 681d4a9c,
# rId3 -> 8681479b) stay put; only which row -- and therefore which
# display text -- sits on B2/B3 changes. Range(...).Hyperlinks.Delete()
# clears every hyperlink on the sheet, so re-add both, in the original
# rId order, right after.
$ov.Range("A1").Hyperlinks.Delete()
$ov.Hyperlinks.Add($ov.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test1/blob/66f9683b0c6f78332286b395d685818097363867/e2e/681d4a9c-e7c4-4811-aa17-cc08b7db6c9f.md", [Type]::Missing, [Type]::Missing, "e2e\681d4a9c-e7c4-4811-aa17-cc08b7db6c9f.md")
$ov.Hyperlinks.Add($ov.Range("B2"), "https://github.com/OpenLocalizationTestOrg/ol-test1/blob/66f9683b0c6f78332286b395d685818097363867/e2e/8681479b-ec15-4532-9f35-a51ee3b9248b.md", [Type]::Missing, [Type]::Missing, "e2e\8681479b-ec15-4532-9f35-a51ee3b9248b.md")

# Columns E/F widened to fit the longer "Ready for handoff" text.
$ov.Columns("E").ColumnWidth = 17.2159881591797
$ov.Columns("F").ColumnWidth = 17.2159881591797

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

$zh.Range("A2").Value = "8681479b-ec15-4532-9f35-a51ee3b9248b.md"
$zh.Range("B2").Value = ".md"
$zh.Range("C2").Value = "Ready for handoff"
$zh.Range("D2").Value = "e2e"
$zh.Range("E2").Value = "ht"
Set-TextValue $zh.Range("F2") "False"
$zh.Range("G2").Value = "8681479b-ec15-4532-9f35-a51ee3b9248b.83253c021a3fcc17a53cec25d51e1b5986801abc.zh-cn.xlf"
$zh.Range("H2").Value = "2017-01-03 05:56:12"
$zh.Range("L2").Value = "0001-01-01 00:00:00"
Set-TextValue $zh.Range("O2") "True"
Set-TextValue $zh.Range("Q2") "False"

$zh.Range("A3").Value = "681d4a9c-e7c4-4811-aa17-cc08b7db6c9f.md"
$zh.Range("B3").Value = ".md"
$zh.Range("C3").Value = "Ready for handoff"
$zh.Range("D3").Value = "e2e"
$zh.Range("E3").Value = "mt"
Set-TextValue $zh.Range("F3") "False"
$zh.Range("G3").Value = "681d4a9c-e7c4-4811-aa17-cc08b7db6c9f.fde05580433d97612bb91a28b095105a51da3047.zh-cn.xlf"
$zh.Range("H3").Value = "2017-01-03 05:57:00"
$zh.Range("L3").Value = "0001-01-01 00:00:00"
Set-TextValue $zh.Range("O3") "True"
Set-TextValue $zh.Range("Q3") "False"

$zh.Range("A1").Hyperlinks.Delete()
$zh.Hyperlinks.Add($zh.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test1/blob/66f9683b0c6f78332286b395d685818097363867/e2e/681d4a9c-e7c4-4811-aa17-cc08b7db6c9f.md", [Type]::Missing, [Type]::Missing, "681d4a9c-e7c4-4811-aa17-cc08b7db6c9f.md")
$zh.Hyperlinks.Add($zh.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test1/blob/66f9683b0c6f78332286b395d685818097363867/e2e/8681479b-ec15-4532-9f35-a51ee3b9248b.md", [Type]::Missing, [Type]::Missing, "8681479b-ec15-4532-9f35-a51ee3b9248b.md")

$zh.Columns("C").ColumnWidth = 17.2159881591797

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

$de.Range("A2").Value = "8681479b-ec15-4532-9f35-a51ee3b9248b.md"
$de.Range("B2").Value = ".md"
$de.Range("C2").Value = "Ready for handoff"
$de.Range("D2").Value = "e2e"
$de.Range("E2").Value = "ht"
Set-TextValue $de.Range("F2") "False"
$de.Range("G2").Value = "8681479b-ec15-4532-9f35-a51ee3b9248b.83253c021a3fcc17a53cec25d51e1b5986801abc.de-de.xlf"
$de.Range("H2").Value = "2017-01-03 05:56:20"
$de.Range("L2").Value = "0001-01-01 00:00:00"
Set-TextValue $de.Range("O2") "True"
Set-TextValue $de.Range("Q2") "False"

$de.Range("A3").Value = "681d4a9c-e7c4-4811-aa17-cc08b7db6c9f.md"
$de.Range("B3").Value = ".md"
$de.Range("C3").Value = "Ready for handoff"
$de.Range("D3").Value = "e2e"
$de.Range("E3").Value = "mt"
Set-TextValue $de.Range("F3") "False"
$de.Range("G3").Value = "681d4a9c-e7c4-4811-aa17-cc08b7db6c9f.fde05580433d97612bb91a28b095105a51da3047.de-de.xlf"
$de.Range("H3").Value = "2017-01-03 05:57:09"
$de.Range("L3").Value = "0001-01-01 00:00:00"
Set-TextValue $de.Range("O3") "True"
Set-TextValue $de.Range("Q3") "False"

$de.Range("A1").Hyperlinks.Delete()
$de.Hyperlinks.Add($de.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test1/blob/66f9683b0c6f78332286b395d685818097363867/e2e/681d4a9c-e7c4-4811-aa17-cc08b7db6c9f.md", [Type]::Missing, [Type]::Missing, "681d4a9c-e7c4-4811-aa17-cc08b7db6c9f.md")
$de.Hyperlinks.Add($de.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test1/blob/66f9683b0c6f78332286b395d685818097363867/e2e/8681479b-ec15-4532-9f35-a51ee3b9248b.md", [Type]::Missing, [Type]::Missing, "8681479b-ec15-4532-9f35-a51ee3b9248b.md")

$de.Columns("C").ColumnWidth = 17.2159881591797
